# Applies the crypto tracker refresh described by the commit:
# "Updated cryptos list on Fri Oct 13 13:33:53 UTC 2023 with GitHub Actions"
#
# Price (column D) and Volume(1h) (column E) figures are refreshed for every
# coin row, and two rows (16/17) swap rank positions (Litecoin <-> WrappedBTC).
#
# All source cells are stored as plain text (coinranking.com scrapes numbers
# like "27.031.68" or "1.10" as literal strings, not locale-aware numerics),
# so each write forces a Text number format before assigning the value and
# then restores the default "Normal" style afterwards -- this keeps Excel from
# auto-coercing dotted strings (e.g. "1.10" or "27.031.68") into numbers while
# leaving cell styling exactly as it was originally (no explicit style index).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @(
    @{ Cell = "D2"; Value = "27.031.68" }
    @{ Cell = "E2"; Value = "  +0.84%  " }
    @{ Cell = "D3"; Value = "1.559.64" }
    @{ Cell = "E3"; Value = "  +0.62%  " }
    @{ Cell = "E4"; Value = "  +0.40%  " }
    @{ Cell = "D5"; Value = "207.43" }
    @{ Cell = "E5"; Value = "  +0.66%  " }
    @{ Cell = "E6"; Value = "  +0.96%  " }
    @{ Cell = "E7"; Value = "  +0.36%  " }
    @{ Cell = "D8"; Value = "21.73" }
    @{ Cell = "E8"; Value = "  +1.49%  " }
    @{ Cell = "E9"; Value = "  +1.18%  " }
    @{ Cell = "E10"; Value = "  +1.54%  " }
    @{ Cell = "D11"; Value = "0.0862" }
    @{ Cell = "E11"; Value = "  +0.79%  " }
    @{ Cell = "E12"; Value = "  +0.51%  " }
    @{ Cell = "D13"; Value = "1.558.69" }
    @{ Cell = "E13"; Value = "  +0.39%  " }
    @{ Cell = "E14"; Value = "  +1.35%  " }
    @{ Cell = "E15"; Value = "  +0.94%  " }
    @{ Cell = "B16"; Value = "Litecoin" }
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc" }
    @{ Cell = "D16"; Value = "62.06" }
    @{ Cell = "E16"; Value = "  +1.36%  " }
    @{ Cell = "B17"; Value = "WrappedBTC" }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc" }
    @{ Cell = "D17"; Value = "27.026.72" }
    @{ Cell = "E17"; Value = "  +0.81%  " }
    @{ Cell = "D18"; Value = "216.26" }
    @{ Cell = "E18"; Value = "  +0.62%  " }
    @{ Cell = "E19"; Value = "  +0.24%  " }
    @{ Cell = "D20"; Value = "7.29" }
    @{ Cell = "E20"; Value = "  +0.48%  " }
    @{ Cell = "E21"; Value = "  +0.42%  " }
    @{ Cell = "E22"; Value = "  -1.19%  " }
    @{ Cell = "D23"; Value = "9.24" }
    @{ Cell = "E23"; Value = "  +2.78%  " }
    @{ Cell = "E24"; Value = "  -1.02%  " }
    @{ Cell = "D25"; Value = "152.89" }
    @{ Cell = "E25"; Value = "  -0.29%  " }
    @{ Cell = "E26"; Value = "  +2.26%  " }
    @{ Cell = "D27"; Value = "14.96" }
    @{ Cell = "E27"; Value = "  +0.22%  " }
    @{ Cell = "E28"; Value = "  +0.39%  " }
    @{ Cell = "E29"; Value = "  +1.46%  " }
    @{ Cell = "E30"; Value = "  +0.58%  " }
    @{ Cell = "D31"; Value = "1.10" }
    @{ Cell = "E31"; Value = "  -0.83%  " }
    @{ Cell = "E32"; Value = "  +1.25%  " }
    @{ Cell = "D33"; Value = "1.404.41" }
    @{ Cell = "E33"; Value = "  +3.95%  " }
    @{ Cell = "E34"; Value = "  +3.33%  " }
    @{ Cell = "E35"; Value = "  +3.33%  " }
    @{ Cell = "D36"; Value = "0.966" }
    @{ Cell = "E36"; Value = "  +3.57%  " }
    @{ Cell = "E37"; Value = "  +0.13%  " }
    @{ Cell = "E38"; Value = "  +1.43%  " }
    @{ Cell = "D39"; Value = "0.524" }
    @{ Cell = "E39"; Value = "  +0.73%  " }
    @{ Cell = "D40"; Value = "0.812" }
    @{ Cell = "E40"; Value = "  +1.32%  " }
    @{ Cell = "E41"; Value = "  +0.39%  " }
    @{ Cell = "E42"; Value = "  -0.28%  " }
    @{ Cell = "E43"; Value = "  +3.48%  " }
    @{ Cell = "E44"; Value = "  -3.19%  " }
    @{ Cell = "D45"; Value = "64.07" }
    @{ Cell = "E45"; Value = "  +1.49%  " }
    @{ Cell = "E46"; Value = "  -0.43%  " }
    @{ Cell = "D47"; Value = "1.694.18" }
    @{ Cell = "E47"; Value = "  +0.42%  " }
    @{ Cell = "D48"; Value = "86.23" }
    @{ Cell = "E48"; Value = "  +0.43%  " }
    @{ Cell = "D49"; Value = "0.0512" }
    @{ Cell = "E49"; Value = "  +0.48%  " }
    @{ Cell = "D50"; Value = "0.0962" }
    @{ Cell = "E50"; Value = "  +1.22%  " }
    @{ Cell = "E51"; Value = "  +0.47%  " }
)

foreach ($update in $cellUpdates) {
    $range = $ws.Range($update.Cell)
    $range.NumberFormat = "@"
    $range.Value = $update.Value
    $range.Style = "Normal"
}
